# Add a "Save" column (H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1 - style matches the other header cells (bold/bordered "style 1")
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Rows with the two highest "sum" (G) values get Save = 1, all others get Save = 0.
$saveValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 1
    18 = 0
    19 = 0
    20 = 0
    21 = 1
    22 = 0
    23 = 0
    24 = 0
    25 = 0
    26 = 0
    27 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
